$d = $word.ActiveDocument

function Replace-Text($find, $replace, $wholeWord = $false) {
    $d.Content.Find.Execute($find, $true, $wholeWord, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "The ants problem - subtitles:" "Tatizo la mchwa - manukuu:"
Replace-Text "The dialogue starts at 40 seconds in so I added 27 seconds to the times as they were - John Argentino" "Mazungumzo huanza kwa sekunde 40 kwa hivyo niliongeza sekunde 27 kwa nyakati kama zilivyokuwa - John Argentino"
Replace-Text "[Music]" "[Muziki]"
Replace-Text "okay so the puzzles I'm going to" "sawa kwa hivyo mafumbo nitaenda"
Replace-Text "challenge you with are two basic" "changamoto uliyonayo ni mbili za msingi"
Replace-Text "versions of a more complicated puzzle" "matoleo ya fumbo ngumu zaidi"
Replace-Text "known as the ants puzzle, which I'm" "inayojulikana kama fumbo la mchwa, ambalo mimi ni"
Replace-Text "probably going to discuss in a different" "pengine kwenda kujadili katika tofauti"
Replace-Text "video. Let me just finish writing down" "video. Ngoja nimalizie kuandika"
Replace-Text "the title and, well, I can even draw a" "kichwa na, vizuri, naweza hata kuchora a"
Replace-Text "little ant right here. okay, let's get" "mchwa mdogo hapa. sawa, tupate"
Replace-Text "started! As I said I'm going to discuss" "imeanza! Kama nilivyosema nitajadili"
Replace-Text "two puzzles in the first puzzle there" "mafumbo mawili katika fumbo la kwanza hapo"
Replace-Text "are two ants on a very high stool: a sort" "ni mchwa wawili kwenye kinyesi cha juu sana: aina"
Replace-Text "of Mountain, flat on the top with two" "ya Mlima, gorofa juu na mbili"
Replace-Text "steep cliffs to both the sides. The flat" "miamba mikali kwa pande zote mbili. Gorofa"
Replace-Text "peak is one meter wide the two ants move" "kilele ni mita moja upana wa mchwa wawili hoja"
Replace-Text "with a velocity, let's call it V, which is" "kwa kasi, tuiite V, ambayo ni"
Replace-Text "the same for both of them and that is" "sawa kwa wote wawili na hiyo ni"
Replace-Text "equal to one centimeter per second. You" "sawa na sentimita moja kwa sekunde. Wewe"
Replace-Text "can decide the direction towards each" "inaweza kuamua mwelekeo kuelekea kila mmoja"
Replace-Text "ant moves if it is right or left and" "mchwa husogea ikiwa ni kulia au kushoto na"
Replace-Text "where exactly to place the two ants on the" "wapi hasa kuweka mchwa wawili kwenye"
Replace-Text "top of the mountain. Your purpose is to" "juu ya mlima. Kusudi lako ni"
Replace-Text "make the time the last ant takes before" "fanya wakati mchwa wa mwisho huchukua hapo awali"
Replace-Text "falling the longest possible. Ants cannot" "kuanguka kwa muda mrefu iwezekanavyo. Mchwa hawawezi"
Replace-Text "be still: they must move to the right or" "tulia: lazima wahamie kulia au"
Replace-Text "to the left but they must move and after" "upande wa kushoto lakini lazima wasogee na baada"
Replace-Text "meeting each other they turn around and" "wakikutana wanageuka na"
Replace-Text "keep moving with the same but opposite" "endelea kusonga na sawa lakini kinyume"
Replace-Text "so again what are the precise positions" "kwa hivyo tena ni nafasi gani sahihi"
Replace-Text "where I should place the two ants in" "ambapo ninapaswa kuwaweka mchwa wawili ndani"
Replace-Text "order to get the longest time before the" "ili kupata muda mrefu zaidi kabla ya"
Replace-Text "last ant falls? The second puzzle is" "chungu mwisho huanguka? Fumbo la pili ni"
Replace-Text "basically the same but now we have three" "kimsingi ni sawa lakini sasa tuna tatu"
Replace-Text "ants instead of two." "mchwa badala ya wawili."
Replace-Text "As before the ants velocity is one" "Kama kabla ya mchwa kasi ni moja"
Replace-Text "centimeter per second, every ant turns" "sentimita kwa sekunde, kila mchwa hugeuka"
Replace-Text "around after meeting another ant and" "karibu baada ya kukutana na mchwa mwingine na"
Replace-Text "the peak is one meter wide. So, what are" "kilele kina upana wa mita moja. Hivyo, ni nini"
Replace-Text "now the precise positions" "sasa nafasi sahihi"
Replace-Text "I should place the three ants in order" "Ninapaswa kuweka mchwa watatu kwa mpangilio"
Replace-Text "to get the longest time before the last" "kupata muda mrefu zaidi kabla ya mwisho"
Replace-Text "ant falls down? I hope you enjoyed this" "chungu huanguka chini? Natumaini ulifurahia hili"
Replace-Text "video do your best and good luck" "video fanya bora na bahati nzuri"
Replace-Text "velocity" "kasi" $true
